# Temporary Codes.xlsx - add a "JAVA" column (D) with generated Java enum-style
# literals to each of the four lookup sheets, and refresh the test-case view
# state (active tab / selections / column widths) to match the authored edit.

$wb = $excel.ActiveWorkbook

$sheetLOINC = $wb.Worksheets.Item("LOINC")            # sheet1.xml
$sheetPandemic = $wb.Worksheets.Item("Pandemic")      # sheet2.xml
$sheetPopGroups = $wb.Worksheets.Item("Population Groups") # sheet3.xml
$sheetTier = $wb.Worksheets.Item("Tier")              # sheet4.xml

# ---------------------------------------------------------------------------
# LOINC sheet (sheet1): rows 2-4, D formula has no SUBSTITUTE (codes have no
# dashes to swap out).
# ---------------------------------------------------------------------------
$ws = $sheetLOINC
$ws.Range("B1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "JAVA"

$ws.Range("D2").Formula = '=A2&"("&CHAR(34)&A2&CHAR(34)&", "&CHAR(34)&B2&CHAR(34)&"),"'
$ws.Range("D3:D4").Formula = '=A3&"("&CHAR(34)&A3&CHAR(34)&", "&CHAR(34)&B3&CHAR(34)&"),"'

$ws.Columns.Item(3).ColumnWidth = 51.333333333333336

$ws.Range("D2:D4").Select()

# ---------------------------------------------------------------------------
# Pandemic sheet (sheet2): rows 2-4, D formula substitutes "-" with "_".
# ---------------------------------------------------------------------------
$ws = $sheetPandemic
$ws.Range("B1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "JAVA"

$ws.Range("D2").Formula = '=SUBSTITUTE(A2, "-", "_")&"("&CHAR(34)&A2&CHAR(34)&", "&CHAR(34)&B2&CHAR(34)&"),"'
$ws.Range("D3:D4").Formula = '=SUBSTITUTE(A3, "-", "_")&"("&CHAR(34)&A3&CHAR(34)&", "&CHAR(34)&B3&CHAR(34)&"),"'

$ws.Columns.Item(3).ColumnWidth = 64.66666666666667

$ws.Range("D2").Select()

# ---------------------------------------------------------------------------
# Population Groups sheet (sheet3): rows 2-29, D formula substitutes "-" with "_".
# ---------------------------------------------------------------------------
$ws = $sheetPopGroups
$ws.Range("B1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "JAVA"

$ws.Range("D2").Formula = '=SUBSTITUTE(A2, "-", "_")&"("&CHAR(34)&A2&CHAR(34)&", "&CHAR(34)&B2&CHAR(34)&"),"'
$ws.Range("D3:D29").Formula = '=SUBSTITUTE(A3, "-", "_")&"("&CHAR(34)&A3&CHAR(34)&", "&CHAR(34)&B3&CHAR(34)&"),"'

$ws.Columns.Item(3).ColumnWidth = 110.16666666666667

$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("D1:D2").Select()
$ws.Range("D1").Activate()

# ---------------------------------------------------------------------------
# Tier sheet (sheet4): rows 2-7, D formula substitutes "-" with "_". This is
# the sheet active when the workbook was last saved.
# ---------------------------------------------------------------------------
$ws = $sheetTier
$ws.Range("B1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "JAVA"

$ws.Range("D2").Formula = '=SUBSTITUTE(A2, "-", "_")&"("&CHAR(34)&A2&CHAR(34)&", "&CHAR(34)&B2&CHAR(34)&"),"'
$ws.Range("D3:D7").Formula = '=SUBSTITUTE(A3, "-", "_")&"("&CHAR(34)&A3&CHAR(34)&", "&CHAR(34)&B3&CHAR(34)&"),"'

$ws.Columns.Item(3).ColumnWidth = 22
$ws.Columns.Item(4).ColumnWidth = 23.5

$ws.Activate()
$ws.Range("D9").Select()
